# temp solve of RWheel
# Set the Fitness column (C2:C12) to a constant value of 4297

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C12").Value = 4297
